$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the recruitment name in F2 (shared string "YTR" -> "RECURIMENT2")
$ws.Range("F2").Value = "RECURIMENT2"

# Move the active selection to G6 (was G17)
$ws.Range("G6").Select()
